# Apply updated cryptos list values (price & 1h volume change) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.917.09"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "3.585.50"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'586.93"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").Value = "'185.00"
$ws.Range("E6").Value = "  -2.39%  "

$ws.Range("D7").Value = "3.575.39"
$ws.Range("E7").Value = "  -1.13%  "

$ws.Range("D8").Value = "'0.622"
$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "'0.217"
$ws.Range("E10").Value = "  +15.11%  "

$ws.Range("E11").Value = "  -1.23%  "

$ws.Range("D12").Value = "'54.22"
$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("D13").Value = "'0.0000321"
$ws.Range("E13").Value = "  +4.60%  "

$ws.Range("E14").Value = "  -1.92%  "

$ws.Range("D15").Value = "4.156.96"
$ws.Range("E15").Value = "  -1.64%  "

$ws.Range("D16").Value = "70.965.23"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").Value = "'19.29"
$ws.Range("E17").Value = "  -2.11%  "

$ws.Range("D18").Value = "3.578.90"
$ws.Range("E18").Value = "  -1.68%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'12.40"
$ws.Range("E19").Value = "  -1.74%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'565.91"
$ws.Range("E20").Value = "  +12.30%  "

$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("E22").Value = "  -4.27%  "

$ws.Range("D23").Value = "'17.63"
$ws.Range("E23").Value = "  -8.86%  "

$ws.Range("E24").Value = "  +4.39%  "

$ws.Range("D25").Value = "'4.99"
$ws.Range("E25").Value = "  +1.17%  "

$ws.Range("D26").Value = "'94.67"
$ws.Range("E26").Value = "  -1.39%  "

$ws.Range("D27").Value = "'11.19"
$ws.Range("E27").Value = "  -3.67%  "

$ws.Range("E28").Value = "  -2.58%  "

$ws.Range("E29").Value = "  -3.45%  "

$ws.Range("D30").Value = "'32.33"
$ws.Range("E30").Value = "  +1.10%  "

$ws.Range("D31").Value = "'7.24"
$ws.Range("E31").Value = "  -6.93%  "

$ws.Range("D32").Value = "'12.28"
$ws.Range("E32").Value = "  -2.73%  "

$ws.Range("D33").Value = "'64.21"
$ws.Range("E33").Value = "  -3.44%  "

$ws.Range("E34").Value = "  -2.54%  "

$ws.Range("D35").Value = "'3.28"
$ws.Range("E35").Value = "  +1.11%  "

$ws.Range("D36").Value = "'550.50"
$ws.Range("E36").Value = "  -4.59%  "

$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").Value = "0.0₃0803"
$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").Value = "'37.52"
$ws.Range("E40").Value = "  -4.23%  "

$ws.Range("D41").Value = "3.462.61"
$ws.Range("E41").Value = "  +7.22%  "

$ws.Range("E42").Value = "  -1.11%  "

$ws.Range("D43").Value = "'3.39"
$ws.Range("E43").Value = "  -2.41%  "

$ws.Range("D44").Value = "'3.11"
$ws.Range("E44").Value = "  -6.53%  "

$ws.Range("D45").Value = "'3.51"
$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("E46").Value = "  -3.30%  "

$ws.Range("E47").Value = "  -2.94%  "

$ws.Range("D48").Value = "'9.33"
$ws.Range("E48").Value = "  -2.95%  "

$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("D50").Value = "'0.997"
$ws.Range("E50").Value = "  -0.45%  "

$ws.Range("D51").Value = "'1.43"
$ws.Range("E51").Value = "  -6.90%  "
